$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Current layout:  Sheet1 = "总计" (summary), Sheet2 = "2022-Q3" (fund table)
# Target layout:   Sheet1 = "总计", Sheet2 = "2022-Q4" (new fund table),
#                   Sheet3 = "2022-Q3" (old fund table, unchanged, moved here)
# ---------------------------------------------------------------------------

$summary = $wb.Worksheets.Item("总计")
$q3 = $wb.Worksheets.Item("2022-Q3")

# Step 1: duplicate the existing "2022-Q3" sheet right after itself so the
# old data survives under its original name on a new tab.
$q3.Copy($null, $q3)
$q3copy = $wb.Worksheets.Item("2022-Q3 (2)")
$q3.Name = "2022-Q4"
$q3copy.Name = "2022-Q3"

$q4 = $q3

# Step 2: wipe the (now renamed) "2022-Q4" sheet and populate it with the
# new quarterly fund-holding data.
$q4.Cells.Clear()

# Columns B (fund code, leading zeros) and D-G (numeric-looking text such as
# "48.20" whose trailing zero must not be dropped) must be stored as text,
# exactly like the source data - force text format before typing values in.
$q4.Range("B2:B6").NumberFormat = "@"
$q4.Range("D2:G6").NumberFormat = "@"

$q4.Range("B1").Value = "基金代码"
$q4.Range("C1").Value = "基金名称"
$q4.Range("D1").Value = "基金规模"
$q4.Range("E1").Value = "股票总仓位"
$q4.Range("F1").Value = "仓位占比"
$q4.Range("G1").Value = "持有市值(亿元)"
$q4.Range("H1").Value = "仓位排名"

# idx, code, name, size, stockPosition, positionRatio, marketValue, rank
$rows = @(
    @(0, "011093", "永赢宏泽一年定期开放灵活配置混合", "14.98", "48.20", "0.39", "0.0584", 9),
    @(1, "519615", "银河君尚灵活配置混合I",           "1.83",  "38.98", "0.99", "0.0181", 5),
    @(2, "006836", "永赢惠泽一年定期开放灵活配置混合", "3.73",  "48.14", "0.39", "0.0145", 10),
    @(3, "519613", "银河君尚灵活配置混合A",           "1.17",  "38.98", "0.99", "0.0116", 5),
    @(4, "519614", "银河君尚灵活配置混合C",           "0.16",  "38.98", "0.99", "0.0016", 5)
)

for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $i + 2
    $row = $rows[$i]
    $q4.Cells.Item($r, 1).Value = $row[0]
    $q4.Cells.Item($r, 2).Value = $row[1]
    $q4.Cells.Item($r, 3).Value = $row[2]
    $q4.Cells.Item($r, 4).Value = $row[3]
    $q4.Cells.Item($r, 5).Value = $row[4]
    $q4.Cells.Item($r, 6).Value = $row[5]
    $q4.Cells.Item($r, 7).Value = $row[6]
    $q4.Cells.Item($r, 8).Value = $row[7]
}

# Now that the text is locked in, drop the explicit "@" number format again
# so the cells end up with the default (General) style, matching the source
# workbook's formatting.
$q4.Range("B2:B6").Style = "Normal"
$q4.Range("D2:G6").Style = "Normal"

# Step 3: re-apply the header / index-column styling used elsewhere in this
# workbook (bold, centered, thin-bordered) by copying formats from the
# summary sheet, which already carries that style.
$summary.Range("B1").Copy()
$q4.Range("B1:H1").PasteSpecial(-4122)   # xlPasteFormats
$summary.Range("A2").Copy()
$q4.Range("A2:A6").PasteSpecial(-4122)   # xlPasteFormats

$excel.CutCopyMode = $false

# Step 4: update the "总计" summary sheet - row 2 now reports 2022-Q4 totals,
# and a new row 3 keeps the original 2022-Q3 totals.
$summary.Range("B2").Value = "2022-Q4"
$summary.Range("C2").Value = 5
$summary.Range("D2").Value = 0.1

$summary.Range("A3").Value = 1
$summary.Range("B3").Value = "2022-Q3"
$summary.Range("C3").Value = 7
$summary.Range("D3").Value = 0.07000000000000001

$summary.Range("A2").Copy()
$summary.Range("A3").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

Write-Host "Workbook updated."
